# Applies the "stoee credentials in EEPROM" commit:
#  1. Bumps the cached datetimeFigureOut placeholder text from 23/01/2019
#     to 24/01/2019 on the slide master and every slide layout.
#  2. Repositions/resizes a handful of shapes on slide 2 (the EEPROM /
#     credentials wiring diagram) and bumps the "34" ellipse label to "38".

$p = $ppt.ActivePresentation
$emuPerPt = 12700

# Shape.Left/Top/Width/Height round-trip through a single-precision (Single)
# COM property and this host truncates (rather than rounds) when it turns
# the point value back into EMU on save. Biasing the EMU value up by half a
# unit before converting to points keeps the stored EMU exact.
function Emu-ToPoints {
    param([double]$emu)
    return ($emu + 0.5) / $emuPerPt
}

function Set-DateText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $txt = $shp.TextFrame.TextRange.Text
            if ($txt -eq "23/01/2019") {
                $shp.TextFrame.TextRange.Text = "24/01/2019"
            }
        }
    }
}

# --- 1. Date placeholder on the slide master ---
Set-DateText $p.SlideMaster.Shapes

# --- 1b. Date placeholder on every slide layout ---
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Set-DateText $layout.Shapes
}

# --- 2. Shape tweaks on slide 2 ---
$s = $p.Slides.Item(2)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)

    switch ($shp.Id) {
        173 {
            # Rectangle 172 ("°C & %H DHT22") — move only
            $shp.Left = Emu-ToPoints 10286361
            $shp.Top  = Emu-ToPoints 3693338
        }
        178 {
            # Ellipse 177 — move + relabel 34 -> 38
            $shp.Left = Emu-ToPoints 9673803
            $shp.Top  = Emu-ToPoints 3730704
            $shp.TextFrame.TextRange.Text = "38"
        }
        179 {
            # Connecteur : en angle 178 — move + resize
            $shp.Left   = Emu-ToPoints 10153636
            $shp.Top    = Emu-ToPoints 3830051
            $shp.Width  = Emu-ToPoints 132725
            $shp.Height = Emu-ToPoints 10673
        }
        182 {
            # Connecteur droit 181 — move only
            $shp.Left = Emu-ToPoints 11525257
            $shp.Top  = Emu-ToPoints 3742400
        }
        183 {
            # Connecteur droit 182 — move only
            $shp.Left = Emu-ToPoints 11517548
            $shp.Top  = Emu-ToPoints 3856925
        }
        224 {
            # Ellipse 223 — move only
            $shp.Left = Emu-ToPoints 8902646
            $shp.Top  = Emu-ToPoints 3599542
        }
        225 {
            # Connecteur : en angle 224 — move + resize
            $shp.Left   = Emu-ToPoints 6927076
            $shp.Top    = Emu-ToPoints 1733992
            $shp.Width  = Emu-ToPoints 2589556
            $shp.Height = Emu-ToPoints 1361583
        }
    }
}
